{"js": "// Sprint 2 report update:\n//  1. \"...members of the group, with the exception of Sebastian...\"\n//     -> \"...members of the group, except for Sebastian...\"\n//  2. \"...single application \u2013  and started adding...\" (double space\n//     before \"and\") -> \"...single application \u2013 and started adding...\"\n//     (single space)\n\nconst body = context.document.body;\n\n// --- Edit 1: \"with the exception of\" -> \"except for\" ------------------\nconst ex1 = body.search(\"with the exception of\", { matchCase: true });\nex1.load(\"text\");\nawait context.sync();\n\nif (ex1.items.length > 0) {\n  ex1.items[0].insertText(\"except for\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Edit 2: drop the extra space before \"and started\" ----------------\n// Search narrowly for \" and start\" (single leading space) so we only\n// touch the run that holds it, then replace it with \"and start\"\n// (no leading space) to collapse the accidental double space that\n// preceded it.\nconst ex2 = body.search(\" and start\", { matchCase: true });\nex2.load(\"text\");\nawait context.sync();\n\nif (ex2.items.length > 0) {\n  ex2.items[0].insertText(\"and start\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Sprint 2 report update:\n#  1. \"...members of the group, with the exception of Sebastian...\"\n#     -> \"...members of the group, except for Sebastian...\"\n#  2. \"...single application -  and started adding...\" (double space\n#     before \"and\") -> \"...single application - and started adding...\"\n#     (single space)\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: \"with the exception of\" -> \"except for\" -------------------\n$rng1 = $d.Content\n$find1 = $rng1.Find\n$find1.ClearFormatting()\n$find1.Text = \"with the exception of\"\n$find1.MatchCase = $true\n$find1.Forward = $true\n$find1.Wrap = 0          # wdFindStop - don't wrap around, single match expected\n$found1 = $find1.Execute()\nif ($found1) {\n    $rng1.Text = \"except for\"\n}\n\n# --- Edit 2: drop the extra space before \"and started\" -----------------\n$rng2 = $d.Content\n$find2 = $rng2.Find\n$find2.ClearFormatting()\n$find2.Text = \" and start\"\n$find2.MatchCase = $true\n$find2.Forward = $true\n$find2.Wrap = 0\n$found2 = $find2.Execute()\nif ($found2) {\n    $rng2.Text = \"and start\"\n}\n"}
